$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source stores every Price/Volume cell as plain text (inline strings),
# never as a number. Assigning a numeric-looking string (e.g. '533.76') via
# .Value normally auto-converts it to a Double, which would change both the
# cell type and the stored text (trailing zeros, float noise, etc).
# Briefly flip each touched Price cell to Text format before the write so
# the literal string is preserved, then restore the default (Normal) style
# afterwards so no stray number-format/style is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '57.935.33'
$ws.Range("E2").Value = '  +1.04%  '
$ws.Range("D3").Value = '3.127.65'
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '533.76'
$ws.Range("E5").Value = '  +1.51%  '
$ws.Range("E6").Value = '  +0.69%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '0.483'
$ws.Range("E8").Value = '  +6.98%  '
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("E10").Value = '  +0.43%  '
$ws.Range("E11").Value = '  +2.86%  '
$ws.Range("D13").Value = '3.666.87'
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("D14").Value = '25.90'
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("E15").Value = '  +1.95%  '
$ws.Range("D16").Value = '58.018.78'
$ws.Range("E16").Value = '  +1.00%  '
$ws.Range("D17").Value = '3.123.95'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").Value = '6.10'
$ws.Range("E18").Value = '  +2.94%  '
$ws.Range("D19").Value = '12.84'
$ws.Range("E19").Value = '  +2.29%  '
$ws.Range("D20").Value = '8.14'
$ws.Range("D21").Value = '375.61'
$ws.Range("E21").Value = '  +7.21%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").Value = '5.73'
$ws.Range("E23").Value = '  -1.37%  '
$ws.Range("D24").Value = '69.57'
$ws.Range("E24").Value = '  +1.74%  '
$ws.Range("D25").Value = '0.508'
$ws.Range("E25").Value = '  +0.93%  '
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("D28").Value = '0.0₃0879'
$ws.Range("E28").Value = '  -1.27%  '
$ws.Range("D29").Value = '7.49'
$ws.Range("E29").Value = '  +1.44%  '
$ws.Range("E30").Value = '  +3.15%  '
$ws.Range("E31").Value = '  -0.43%  '
$ws.Range("E32").Value = '  +3.14%  '
$ws.Range("D33").Value = '5.16'
$ws.Range("E33").Value = '  +2.56%  '
$ws.Range("E34").Value = '  +2.86%  '
$ws.Range("E35").Value = '  +0.51%  '
$ws.Range("E36").Value = '  +1.60%  '
$ws.Range("E37").Value = '  +2.99%  '
$ws.Range("D38").Value = '25.61'
$ws.Range("E38").Value = '  -3.11%  '
$ws.Range("D39").Value = '1.65'
$ws.Range("E39").Value = '  +4.23%  '
$ws.Range("E40").Value = '  +2.23%  '
$ws.Range("D41").Value = '2.550.72'
$ws.Range("E41").Value = '  +5.85%  '
$ws.Range("D42").Value = '4.12'
$ws.Range("E42").Value = '  +0.63%  '
$ws.Range("E43").Value = '  +0.57%  '
$ws.Range("E44").Value = '  +3.18%  '
$ws.Range("D45").Value = '0.0271'
$ws.Range("E45").Value = '  +2.49%  '
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").Value = '0.985'
$ws.Range("E47").Value = '  +1.22%  '
$ws.Range("E48").Value = '  +2.73%  '
$ws.Range("D49").Value = '19.82'
$ws.Range("E49").Value = '  +0.32%  '
$ws.Range("D50").Value = '0.749'
$ws.Range("E50").Value = '  -2.17%  '
$ws.Range("E51").Value = '  +3.26%  '

# Restore default styling on the Price cells we touched.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
